$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.679.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.01%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.367.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.46%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.81%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.50%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.20%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.363.71'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.104'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.05%  '

$ws.Range('E11').Value = '  -1.08%  '

$ws.Range('E12').Value = '  +0.23%  '

$ws.Range('E13').Value = '  -0.01%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.86%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.798.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.41%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000164'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.65%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '59.571.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.38%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.362.75'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.08%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('E22').Value = '  +0.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.43%  '

$ws.Range('E24').Value = '  +0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.29%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '558.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.68%  '

$ws.Range('E29').Value = '  -3.34%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0918'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.88%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.02%  '

$ws.Range('E32').Value = '  -2.85%  '

$ws.Range('E33').Value = '  -3.04%  '

$ws.Range('E34').Value = '  -0.98%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.52%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.65%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '152.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.365'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.17%  '

$ws.Range('E39').Value = '  -0.95%  '

$ws.Range('E40').Value = '  -0.14%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.95'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.47%  '

$ws.Range('E42').Value = '  -0.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.62%  '

$ws.Range('E45').Value = '  +3.10%  '

$ws.Range('E46').Value = '  +7.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '138.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.27%  '

$ws.Range('E48').Value = '  +0.65%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.583'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0499'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.75%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.83%  '
